# Generate Report for Handoff
# A new handoff was generated for 3b655d02-e6b9-4377-b28c-9ae3e60d3ae9.md,
# so its "Latest Handoff Date(time)" is refreshed on the Overview sheet and
# on each per-language detail sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D4").Value = "2016-03-18 16:29:33"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-18 16:29:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-18 16:29:33"
